$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 'https://acuiti-preview.abbviepro.com/bin/public/abbvie-commons/basic-login'
$ws.Range("B6").Value = 302
$ws.Range("A7").Value = 'https://cag.abbvie.com:9999/bf/16a183f6-c871-4082-850b-a1f7a2ecd0b1?type=js3&sn=v_4_srv_-2D63_sn_C7A8HPT2NH63CA53SEI0IFGTJVI820KJ&svrid=-63&flavor=cors&vi=QCSNFVHKSMDFCFCMHAHRUKLQJVAADNIB-0&modifiedSince=1713766999397&rf=https%3A%2F%2Facuiti-preview.abbviepro.com%2Fbasic-login.html&bp=3&app=c68a93205a30c351&crc=4207847994&en=oao3vfhf&end=1'
$ws.Range("A8").Value = 'https://cag.abbvie.com:9999/bf/16a183f6-c871-4082-850b-a1f7a2ecd0b1?type=js3&sn=v_4_srv_-2D63_sn_C7A8HPT2NH63CA53SEI0IFGTJVI820KJ&svrid=-63&flavor=cors&vi=QCSNFVHKSMDFCFCMHAHRUKLQJVAADNIB-0&modifiedSince=1713766999397&rf=https%3A%2F%2Facuiti-preview.abbviepro.com%2Fbasic-login.html&bp=3&app=c68a93205a30c351&crc=4263852186&en=oao3vfhf&end=1'
$ws.Range("A9").Value = 'https://cag.abbvie.com:9999/bf/16a183f6-c871-4082-850b-a1f7a2ecd0b1?type=js3&sn=v_4_srv_8_sn_C7A8HPT2NH63CA53SEI0IFGTJVI820KJ_app-3Ac68a93205a30c351_1_ol_0_perc_100000_mul_1&svrid=8&flavor=cors&vi=QCSNFVHKSMDFCFCMHAHRUKLQJVAADNIB-0&modifiedSince=1720290200630&rf=https%3A%2F%2Facuiti-preview.abbviepro.com%2Fbasic-login.html&bp=3&app=c68a93205a30c351&crc=2911120879&en=oao3vfhf&end=1'
$ws.Range("B9").Value = 200
$ws.Range("A17").Value = 'https://acuiti-preview.abbviepro.com/etc.clientlibs/abbvie-pro/clientlibs/publish-header.min.css'
$ws.Range("A18").Value = 'https://acuiti-preview.abbviepro.com/etc.clientlibs/abbvie-pro/clientlibs/publish-header.min.js'
$ws.Range("A30").Value = 'https://acuiti-preview.abbviepro.com/content/abbvie-pro-acuiti/index/jcr:content/body/column_control/par1-100col/image_extension_232433702/item_1.coreimg.png/1662647235751-logo-final_1.png'
$ws.Range("A31").Value = 'https://acuiti-preview.abbviepro.com/content/abbvie-pro-acuiti/index/jcr:content/body/column_control/par1-100col/image_extension/item_1.coreimg.png/1662647235867-Vector.png'
$ws.Range("A32").Value = 'https://cag.abbvie.com:9999/bf/16a183f6-c871-4082-850b-a1f7a2ecd0b1?type=js3&sn=v_4_srv_8_sn_C7A8HPT2NH63CA53SEI0IFGTJVI820KJ_app-3Ac68a93205a30c351_1_ol_0_perc_100000_mul_1&svrid=8&flavor=cors&vi=QCSNFVHKSMDFCFCMHAHRUKLQJVAADNIB-0&modifiedSince=1720290200630&rf=https%3A%2F%2Facuiti-preview.abbviepro.com%2F&bp=3&app=c68a93205a30c351&crc=3875169855&en=oao3vfhf&end=1'
$ws.Range("A33").Value = 'https://cag.abbvie.com:9999/bf/16a183f6-c871-4082-850b-a1f7a2ecd0b1?type=js3&sn=v_4_srv_8_sn_C7A8HPT2NH63CA53SEI0IFGTJVI820KJ_app-3Ac68a93205a30c351_1_ol_0_perc_100000_mul_1&svrid=8&flavor=cors&vi=QCSNFVHKSMDFCFCMHAHRUKLQJVAADNIB-0&modifiedSince=1720290200630&rf=https%3A%2F%2Facuiti-preview.abbviepro.com%2F&bp=3&app=c68a93205a30c351&crc=3875169855&en=oao3vfhf&end=1'
$ws.Range("A34").Value = 'https://cag.abbvie.com:9999/bf/16a183f6-c871-4082-850b-a1f7a2ecd0b1?type=js3&sn=v_4_srv_8_sn_C7A8HPT2NH63CA53SEI0IFGTJVI820KJ_app-3Ac68a93205a30c351_1_ol_0_perc_100000_mul_1&svrid=8&flavor=cors&vi=QCSNFVHKSMDFCFCMHAHRUKLQJVAADNIB-0&modifiedSince=1720290200630&rf=https%3A%2F%2Facuiti-preview.abbviepro.com%2F&bp=3&app=c68a93205a30c351&crc=4254139858&en=oao3vfhf&end=1'
$ws.Range("A35").Value = 'https://cag.abbvie.com:9999/bf/16a183f6-c871-4082-850b-a1f7a2ecd0b1?type=js3&sn=v_4_srv_8_sn_C7A8HPT2NH63CA53SEI0IFGTJVI820KJ_app-3Ac68a93205a30c351_1_ol_0_perc_100000_mul_1&svrid=8&flavor=cors&vi=QCSNFVHKSMDFCFCMHAHRUKLQJVAADNIB-0&modifiedSince=1720290200630&rf=https%3A%2F%2Facuiti-preview.abbviepro.com%2F&bp=3&app=c68a93205a30c351&crc=4254139858&en=oao3vfhf&end=1'
